# CryCompanywiseStockReport_1.xlsx update
#
# The report lists, per company, stock-item rows with columns:
#   A=Sr#  B=ItemCode  C=ItemName  D=Rate  E=MRP  F=Qty  G=Value(=D*Qty)
# followed by a "Sub Total:" row (B = sum of G for the block) and,
# at the very end, an overall "Sub Total:" / "Grand Total:" pair.
#
# This edit re-syncs a batch of rows against a refreshed stock pull:
#   - several adjacent item-row pairs had their ItemCode/MRP/Qty/Value
#     (B/E/F/G) swapped between the two rows (the underlying items kept
#     their Sr#, Name and Rate — only which code/qty landed on which
#     row changed);
#   - several single item rows had their on-hand Qty (F) reduced, with
#     Value (G = D * F) recomputed to match;
#   - the "Sub Total:" rows for the affected blocks, and the workbook's
#     overall Sub Total / Grand Total, were updated to the new sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($RowA, $RowB)
    foreach ($col in @("B", "E", "F", "G")) {
        $addrA = "$col$RowA"
        $addrB = "$col$RowB"
        $valA = $ws.Range($addrA).Value()
        $valB = $ws.Range($addrB).Value()
        $ws.Range($addrA).Value = $valB
        $ws.Range($addrB).Value = $valA
    }
}

# ---- Row pairs whose B/E/F/G (ItemCode/MRP/Qty/Value) got swapped ----
$swapPairs = @(
    @(112, 113),
    @(127, 128),
    @(219, 220),
    @(227, 228),
    @(229, 230),
    @(366, 367),
    @(375, 376),
    @(382, 383),
    @(385, 386),
    @(442, 443),
    @(463, 464)
)

foreach ($pair in $swapPairs) {
    Swap-RowData $pair[0] $pair[1]
}

# ---- Single rows whose on-hand Qty (F) dropped, with Value (G) recomputed ----
$qtyUpdates = @{
    146 = @{ F = 22;   G = 1852.18 }
    174 = @{ F = 9;    G = 2733.75 }
    255 = @{ F = 560;  G = 95944.8 }
    256 = @{ F = 285;  G = 43083.45 }
    291 = @{ F = 112;  G = 4817.12 }
    293 = @{ F = 37;   G = 2601.84 }
    326 = @{ F = 63;   G = 1873.62 }
    328 = @{ F = 42;   G = 1562.82 }
    599 = @{ F = 1638; G = 267174.18 }
    602 = @{ F = 330;  G = 47734.5 }
}

foreach ($row in $qtyUpdates.Keys) {
    $ws.Range("F$row").Value = $qtyUpdates[$row].F
    $ws.Range("G$row").Value = $qtyUpdates[$row].G
}

# ---- "Sub Total:" rows recomputed for the blocks whose Qty/Value changed ----
$subTotals = @{
    147 = 14495.96
    175 = 27971.49
    260 = 189743.97
    304 = 174958.37
    330 = 27838.78
    606 = 431733.43
}

foreach ($row in $subTotals.Keys) {
    $ws.Range("B$row").Value = $subTotals[$row]
}

# ---- Workbook-wide "Sub Total:" / "Grand Total:" rows ----
$ws.Range("B619").Value = 1749545.41
$ws.Range("B620").Value = 1749545.41
